# Updated cryptos list - refresh Price/Volume(1h) values, and update two
# rows where the underlying coin ranking changed position/identity.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where both D (Price) and E (Volume(1h)) change.
# D values are forced to Text format first so Excel doesn't coerce
# numeric-looking strings (e.g. "0.482", "6.13") into floating point
# numbers, and so "thousand-dot" formatted prices (e.g. "2.969.05")
# stay exactly as authored.
$priceVolumeUpdates = @(
    @{ Row = 2;  D = "60.700.90";  E = "  -1.95%  " },
    @{ Row = 3;  D = "2.969.05";   E = "  -1.63%  " },
    @{ Row = 4;  D = "0.999";      E = "  -0.07%  " },
    @{ Row = 5;  D = "521.63";     E = "  -0.88%  " },
    @{ Row = 6;  D = "130.13";     E = "  +1.65%  " },
    @{ Row = 8;  D = "2.964.78";   E = "  -1.63%  " },
    @{ Row = 9;  D = "0.482";      E = "  -1.69%  " },
    @{ Row = 10; D = "6.13";       E = "  +2.51%  " },
    @{ Row = 11; D = "0.146";      E = "  -1.01%  " },
    @{ Row = 12; D = "0.435";      E = "  -1.11%  " },
    @{ Row = 14; D = "32.86";      E = "  -0.48%  " },
    @{ Row = 15; D = "3.442.78";   E = "  -1.54%  " },
    @{ Row = 17; D = "60.692.25";  E = "  -1.87%  " },
    @{ Row = 18; D = "2.971.51";   E = "  -1.43%  " },
    @{ Row = 19; D = "6.44";       E = "  +0.24%  " },
    @{ Row = 20; D = "454.21";     E = "  -3.13%  " },
    @{ Row = 22; D = "0.665";      E = "  -1.85%  " },
    @{ Row = 24; D = "77.61";      E = "  -0.12%  " },
    @{ Row = 25; D = "11.65";      E = "  +0.01%  " },
    @{ Row = 26; D = "0.999";      E = "  -0.07%  " },
    @{ Row = 27; D = "2.61";       E = "  +0.24%  " },
    @{ Row = 28; D = "7.65";       E = "  -2.23%  " },
    @{ Row = 29; D = "1.00";       E = "  +0.19%  " },
    @{ Row = 30; D = "25.01";      E = "  -0.60%  " },
    @{ Row = 33; D = "55.03";      E = "  -1.50%  " },
    @{ Row = 36; D = "5.70";       E = "  -0.94%  " },
    @{ Row = 37; D = "447.55";     E = "  -2.70%  " },
    @{ Row = 38; D = "3.154.47";   E = "  +4.30%  " },
    @{ Row = 39; D = "0.0773";     E = "  +0.44%  " },
    @{ Row = 43; D = "2.39";       E = "  -2.47%  " },
    @{ Row = 46; D = "25.09";      E = "  +6.84%  " },
    @{ Row = 47; D = "118.22";     E = "  +0.59%  " },
    @{ Row = 49; D = "1.92";       E = "  -1.16%  " },
    @{ Row = 50; D = "0.0₃0501";   E = "  -3.47%  " }
)

foreach ($u in $priceVolumeUpdates) {
    $dCell = $ws.Range("D$($u.Row)")
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
    $dCell.Style = "Normal"
    $ws.Range("E$($u.Row)").Value = $u.E
}

# Rows where only E (Volume(1h)) changes
$volumeOnlyUpdates = @(
    @{ Row = 13; E = "  -0.69%  " },
    @{ Row = 16; E = "  +0.47%  " },
    @{ Row = 21; E = "  +0.54%  " },
    @{ Row = 23; E = "  -1.16%  " },
    @{ Row = 31; E = "  +4.10%  " },
    @{ Row = 32; E = "  +1.20%  " },
    @{ Row = 40; E = "  -1.54%  " },
    @{ Row = 41; E = "  +4.60%  " },
    @{ Row = 42; E = "  +1.15%  " },
    @{ Row = 45; E = "  -0.11%  " },
    @{ Row = 48; E = "  +1.09%  " }
)

foreach ($u in $volumeOnlyUpdates) {
    $ws.Range("E$($u.Row)").Value = $u.E
}

# Row 34 and 35 swap identity: NEARProtocol (row34) <-> Stacks (row35)
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$d34 = $ws.Range("D34")
$d34.NumberFormat = "@"
$d34.Value = "2.24"
$d34.Style = "Normal"
$ws.Range("E34").Value = "  -2.60%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$d35 = $ws.Range("D35")
$d35.NumberFormat = "@"
$d35.Value = "5.33"
$d35.Style = "Normal"
$ws.Range("E35").Value = "  +5.37%  "

# Row 51: CoreDAO replaced by BitgetToken
$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$d51 = $ws.Range("D51")
$d51.NumberFormat = "@"
$d51.Value = "1.23"
$d51.Style = "Normal"
$ws.Range("E51").Value = "  +7.82%  "
